# Timesheet sign-off: the supervisor (Prakruti Sinha) has now signed off
# on the week's hours. Fill in:
#   - the Supervisor Name field (G6)
#   - the Supervisor Signature block (initials in A27, signed date in D27)
# Formatting for the newly-filled cells is picked up from the matching,
# already-filled Employee sign-off cells (A25 initials / D25 date) so the
# sheet stays visually consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

$ws.Range("G6").Value = "Prakruti Sinha"

$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41682

$ws.Range("A25").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D25").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
